$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1200.22
$wsSummary.Range("B4").Value = 0.23
$wsSummary.Range("B5").Value = 0.04
$wsSummary.Range("B6").Value = 126
$wsSummary.Range("B8").Value = 48
$wsSummary.Range("B9").Value = 44.44

# --- Strategy Status sheet ---
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 100.22
$wsStatus.Range("D4").Value = 126
$wsStatus.Range("E4").Value = 0.23
$wsStatus.Range("F4").Value = 0.22
$wsStatus.Range("G4").Value = 44.44

# --- All Trades & MarketMaking sheets: append new trade row 127 ---
$sheetNames = @("All Trades", "MarketMaking")
foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("A127").Value = 126
    $ws.Range("B127").NumberFormat = "@"
    $ws.Range("B127").Value = "2026-02-17"
    $ws.Range("C127").Value = "09:32:48"
    $ws.Range("D127").Value = "MarketMaking"
    $ws.Range("E127").Value = "UP"
    $ws.Range("F127").Value = 0.68
    $ws.Range("G127").Value = 0.59
    $ws.Range("H127").Value = "CLOSED"
    $ws.Range("I127").Value = -13.2353
    $ws.Range("J127").Value = -0.09
    $ws.Range("K127").Value = 100.22
    $ws.Range("L127").Value = 0
    $ws.Range("M127").Value = 0
    $ws.Range("N127").Value = 0.6
    $ws.Range("O127").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P127").Value = "early_exit"
    $ws.Range("Q127").Value = 0.1
}
